$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings remain text (matches source formatting,
# e.g. trailing zeros like "0.500"), by pre-setting the column to Text format.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.668.49'
$ws.Range("E2").Value = '  -0.71%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.587.63'
$ws.Range("E3").Value = '  -2.54%  '
$ws.Range("E4").Value = '  +0.50%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '206.96'
$ws.Range("E5").Value = '  -2.07%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.500'
$ws.Range("E6").Value = '  -3.39%  '
$ws.Range("E7").Value = '  +0.51%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '22.20'
$ws.Range("E8").Value = '  -4.66%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.252'
$ws.Range("E9").Value = '  -2.01%  '
$ws.Range("E10").Value = '  -2.74%  '
$ws.Range("E11").Value = '  -1.37%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.813.62'
$ws.Range("E12").Value = '  -2.51%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.571.22'
$ws.Range("E13").Value = '  -3.24%  '
$ws.Range("E14").Value = '  -4.10%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.529'
$ws.Range("E15").Value = '  -4.69%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '27.645.73'
$ws.Range("E16").Value = '  -0.83%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.43'
$ws.Range("E17").Value = '  -2.37%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '219.05'
$ws.Range("E18").Value = '  -4.20%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0₃0696'
$ws.Range("E19").Value = '  -2.95%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.30'
$ws.Range("E20").Value = '  -4.12%  '
$ws.Range("E21").Value = '  +0.49%  '
$ws.Range("E22").Value = '  -4.93%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.59'
$ws.Range("E23").Value = '  -3.48%  '
$ws.Range("E24").Value = '  -4.05%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '153.41'
$ws.Range("E25").Value = '  -1.30%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.83'
$ws.Range("E26").Value = '  -1.25%  '
$ws.Range("E27").Value = '  +0.51%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.10'
$ws.Range("E28").Value = '  -2.33%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.105'
$ws.Range("E29").Value = '  -4.93%  '
$ws.Range("E30").Value = '  -2.35%  '
$ws.Range("E31").Value = '  -2.77%  '
$ws.Range("E32").Value = '  -5.52%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.369.20'
$ws.Range("E33").Value = '  -3.35%  '
$ws.Range("E34").Value = '  -5.69%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.53'
$ws.Range("E35").Value = '  -5.14%  '
$ws.Range("E36").Value = '  -2.56%  '
$ws.Range("E37").Value = '  -0.96%  '
$ws.Range("E38").Value = '  -1.36%  '
$ws.Range("E39").Value = '  -3.33%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.823'
$ws.Range("E40").Value = '  -3.51%  '
$ws.Range("E41").Value = '  +0.52%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.974'
$ws.Range("E42").Value = '  -3.03%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '64.03'
$ws.Range("E43").Value = '  -2.67%  '
$ws.Range("E44").Value = '  +2.49%  '
$ws.Range("E45").Value = '  -4.31%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.724.83'
$ws.Range("E46").Value = '  -2.52%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.71'
$ws.Range("E47").Value = '  -5.77%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '87.53'
$ws.Range("E48").Value = '  -1.29%  '
$ws.Range("E49").Value = '  -1.42%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0965'
$ws.Range("E50").Value = '  -4.68%  '
$ws.Range("E51").Value = '  -1.49%  '
